# Add a new worksheet named "Léo Mastréo" and populate it with a small
# timestamp / status log, matching the uploaded services.xlsx workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add()
$ws.Name = "Léo Mastréo"

# Header row
$ws.Range("A1").Value = "Timestamp"
$ws.Range("B1").Value = "Status"

# Data rows - timestamps are stored as plain text, not Excel dates.
$ws.Range("A2").Value = "2024-07-26T06:03:55.753Z"
$ws.Range("B2").Value = "en service"

$ws.Range("A3").Value = "2024-07-26T06:05:50.261Z"
$ws.Range("B3").Value = "hors service"

$ws.Range("A4").Value = "2024-07-26T06:05:58.286Z"
$ws.Range("B4").Value = "en service"

$ws.Range("A5").Value = "2024-07-26T06:06:30.220Z"
$ws.Range("B5").Value = "hors service"

# Page setup to mirror the source workbook (portrait, fit to one page).
$ws.PageSetup.Orientation = 1
$ws.PageSetup.Zoom = 100
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.FirstPageNumber = 1
$ws.PageSetup.Copies = 1
